# edit.ps1
# Applies the weekend TI (Today's Injuries) roster update:
#  - Removes the SAC / De'Aaron Fox row (team played earlier, no longer listed)
#  - Moves the UTA / Lauri Markkanen row up into its sorted position (row 8)
#    and refreshes his status (Out -> Day-To-Day) and M-1/M-2 numbers
#  - Refreshes rolling stat columns (5M/15M/Saison/GP/age buckets/last-5
#    games/trend columns) for the players that shifted position
#  - Updates Nic Claxton's injury status and Khris Middleton's stat line

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Capture the Lauri Markkanen row (currently row 17) before we shuffle rows ---
$markkanenRow = @()
for ($c = 1; $c -le 38; $c++) {
    $markkanenRow += ,$ws.Cells.Item(17, $c).Value2
}

# --- Remove the SAC / De'Aaron Fox row; everything below shifts up one row ---
$ws.Rows.Item(2).Delete()

# --- Re-insert a row at position 8 (after CHI / Nikola Vucevic) to host Markkanen ---
$ws.Rows.Item(8).Insert()
for ($c = 1; $c -le 38; $c++) {
    $ws.Cells.Item(8, $c).Value2 = $markkanenRow[$c - 1]
}

# --- Drop the now-empty duplicate row that used to hold Markkanen's data ---
$ws.Rows.Item(17).Delete()


# Row 2: MIL / Damian Lillard (refreshed stats)
$ws.Cells.Item(2, 1).Value2 = 'MIL'
$ws.Cells.Item(2, 2).Value2 = 'Damian Lillard'
$ws.Cells.Item(2, 3).Value2 = $null
$ws.Cells.Item(2, 4).Value2 = 'G'
$ws.Cells.Item(2, 5).Value2 = 'O'
$ws.Cells.Item(2, 6).Value2 = 40.8
$ws.Cells.Item(2, 7).Value2 = 35
$ws.Cells.Item(2, 8).Value2 = 35.6
$ws.Cells.Item(2, 9).Value2 = 13
$ws.Cells.Item(2, 10).Value2 = 2
$ws.Cells.Item(2, 11).Value2 = 3
$ws.Cells.Item(2, 12).Value2 = 1
$ws.Cells.Item(2, 13).Value2 = 7
$ws.Cells.Item(2, 14).Value2 = 47
$ws.Cells.Item(2, 15).Value2 = 57
$ws.Cells.Item(2, 16).Value2 = 27
$ws.Cells.Item(2, 17).Value2 = 13
$ws.Cells.Item(2, 18).Value2 = 60
$ws.Cells.Item(2, 19).Value2 = 'vs'
$ws.Cells.Item(2, 20).Value2 = 1.9
$ws.Cells.Item(2, 21).Value2 = -1.4
$ws.Cells.Item(2, 22).Value2 = 9
$ws.Cells.Item(2, 23).Value2 = 'vs'
$ws.Cells.Item(2, 24).Value2 = 'BKN'
$ws.Cells.Item(2, 25).Value2 = '@'
$ws.Cells.Item(2, 26).Value2 = 11
$ws.Cells.Item(2, 27).Value2 = '@'
$ws.Cells.Item(2, 28).Value2 = 33
$ws.Cells.Item(2, 29).Value2 = '-'
$ws.Cells.Item(2, 30).Value2 = '-'
$ws.Cells.Item(2, 31).Value2 = 'vs'
$ws.Cells.Item(2, 32).Value2 = 'OKC'
$ws.Cells.Item(2, 33).Value2 = 'vs'
$ws.Cells.Item(2, 34).Value2 = 'LAL'
$ws.Cells.Item(2, 35).Value2 = '@'
$ws.Cells.Item(2, 36).Value2 = 'NOP'
$ws.Cells.Item(2, 37).Value2 = '@'
$ws.Cells.Item(2, 38).Value2 = 'ATL'

# Row 3: PHX / Devin Booker (refreshed stats)
$ws.Cells.Item(3, 1).Value2 = 'PHX'
$ws.Cells.Item(3, 2).Value2 = 'Devin Booker'
$ws.Cells.Item(3, 3).Value2 = $null
$ws.Cells.Item(3, 4).Value2 = 'G'
$ws.Cells.Item(3, 5).Value2 = 'O'
$ws.Cells.Item(3, 6).Value2 = 34.2
$ws.Cells.Item(3, 7).Value2 = 35
$ws.Cells.Item(3, 8).Value2 = 40.7
$ws.Cells.Item(3, 9).Value2 = 10
$ws.Cells.Item(3, 10).Value2 = 0
$ws.Cells.Item(3, 11).Value2 = 0
$ws.Cells.Item(3, 12).Value2 = 8
$ws.Cells.Item(3, 13).Value2 = 2
$ws.Cells.Item(3, 14).Value2 = 36
$ws.Cells.Item(3, 15).Value2 = 33
$ws.Cells.Item(3, 16).Value2 = 34
$ws.Cells.Item(3, 17).Value2 = 30
$ws.Cells.Item(3, 18).Value2 = 38
$ws.Cells.Item(3, 19).Value2 = 'vs'
$ws.Cells.Item(3, 20).Value2 = -1.7
$ws.Cells.Item(3, 21).Value2 = -12.4
$ws.Cells.Item(3, 22).Value2 = 9
$ws.Cells.Item(3, 23).Value2 = 'vs'
$ws.Cells.Item(3, 24).Value2 = 'ATL'
$ws.Cells.Item(3, 25).Value2 = '@'
$ws.Cells.Item(3, 26).Value2 = 35
$ws.Cells.Item(3, 27).Value2 = '-'
$ws.Cells.Item(3, 28).Value2 = '-'
$ws.Cells.Item(3, 29).Value2 = '-'
$ws.Cells.Item(3, 30).Value2 = '-'
$ws.Cells.Item(3, 31).Value2 = '@'
$ws.Cells.Item(3, 32).Value2 = 'SAS'
$ws.Cells.Item(3, 33).Value2 = '@'
$ws.Cells.Item(3, 34).Value2 = 'SAS'
$ws.Cells.Item(3, 35).Value2 = '@'
$ws.Cells.Item(3, 36).Value2 = 'DEN'
$ws.Cells.Item(3, 37).Value2 = '@'
$ws.Cells.Item(3, 38).Value2 = 'OKC'

# Row 4: ATL / Dejounte Murray (unchanged, re-asserted)
$ws.Cells.Item(4, 1).Value2 = 'ATL'
$ws.Cells.Item(4, 2).Value2 = 'Dejounte Murray'
$ws.Cells.Item(4, 3).Value2 = $null
$ws.Cells.Item(4, 4).Value2 = 'G'
$ws.Cells.Item(4, 5).Value2 = $null
$ws.Cells.Item(4, 6).Value2 = 33.4
$ws.Cells.Item(4, 7).Value2 = 33.5
$ws.Cells.Item(4, 8).Value2 = 31.5
$ws.Cells.Item(4, 9).Value2 = 13
$ws.Cells.Item(4, 10).Value2 = 2
$ws.Cells.Item(4, 11).Value2 = 2
$ws.Cells.Item(4, 12).Value2 = 3
$ws.Cells.Item(4, 13).Value2 = 6
$ws.Cells.Item(4, 14).Value2 = 16
$ws.Cells.Item(4, 15).Value2 = 33
$ws.Cells.Item(4, 16).Value2 = 43
$ws.Cells.Item(4, 17).Value2 = 48
$ws.Cells.Item(4, 18).Value2 = 27
$ws.Cells.Item(4, 19).Value2 = '@'
$ws.Cells.Item(4, 20).Value2 = 0
$ws.Cells.Item(4, 21).Value2 = $null
$ws.Cells.Item(4, 22).Value2 = $null
$ws.Cells.Item(4, 23).Value2 = '@'
$ws.Cells.Item(4, 24).Value2 = 'PHX'
$ws.Cells.Item(4, 25).Value2 = 'vs'
$ws.Cells.Item(4, 26).Value2 = 32
$ws.Cells.Item(4, 27).Value2 = '-'
$ws.Cells.Item(4, 28).Value2 = '-'
$ws.Cells.Item(4, 29).Value2 = '-'
$ws.Cells.Item(4, 30).Value2 = '-'
$ws.Cells.Item(4, 31).Value2 = 'vs'
$ws.Cells.Item(4, 32).Value2 = 'CHA'
$ws.Cells.Item(4, 33).Value2 = 'vs'
$ws.Cells.Item(4, 34).Value2 = 'BOS'
$ws.Cells.Item(4, 35).Value2 = 'vs'
$ws.Cells.Item(4, 36).Value2 = 'POR'
$ws.Cells.Item(4, 37).Value2 = 'vs'
$ws.Cells.Item(4, 38).Value2 = 'BOS'

# Row 5: PHX / Bradley Beal (refreshed stats)
$ws.Cells.Item(5, 1).Value2 = 'PHX'
$ws.Cells.Item(5, 2).Value2 = 'Bradley Beal'
$ws.Cells.Item(5, 3).Value2 = $null
$ws.Cells.Item(5, 4).Value2 = 'G'
$ws.Cells.Item(5, 5).Value2 = 'O'
$ws.Cells.Item(5, 6).Value2 = 30.6
$ws.Cells.Item(5, 7).Value2 = 31.7
$ws.Cells.Item(5, 8).Value2 = 27
$ws.Cells.Item(5, 9).Value2 = 10
$ws.Cells.Item(5, 10).Value2 = 2
$ws.Cells.Item(5, 11).Value2 = 2
$ws.Cells.Item(5, 12).Value2 = 4
$ws.Cells.Item(5, 13).Value2 = 2
$ws.Cells.Item(5, 14).Value2 = 8
$ws.Cells.Item(5, 15).Value2 = 48
$ws.Cells.Item(5, 16).Value2 = 28
$ws.Cells.Item(5, 17).Value2 = 34
$ws.Cells.Item(5, 18).Value2 = 35
$ws.Cells.Item(5, 19).Value2 = 'vs'
$ws.Cells.Item(5, 20).Value2 = -1.6
$ws.Cells.Item(5, 21).Value2 = 0.2
$ws.Cells.Item(5, 22).Value2 = 7
$ws.Cells.Item(5, 23).Value2 = 'vs'
$ws.Cells.Item(5, 24).Value2 = 'ATL'
$ws.Cells.Item(5, 25).Value2 = '@'
$ws.Cells.Item(5, 26).Value2 = 15
$ws.Cells.Item(5, 27).Value2 = '-'
$ws.Cells.Item(5, 28).Value2 = '-'
$ws.Cells.Item(5, 29).Value2 = '-'
$ws.Cells.Item(5, 30).Value2 = '-'
$ws.Cells.Item(5, 31).Value2 = '@'
$ws.Cells.Item(5, 32).Value2 = 'SAS'
$ws.Cells.Item(5, 33).Value2 = '@'
$ws.Cells.Item(5, 34).Value2 = 'SAS'
$ws.Cells.Item(5, 35).Value2 = '@'
$ws.Cells.Item(5, 36).Value2 = 'DEN'
$ws.Cells.Item(5, 37).Value2 = '@'
$ws.Cells.Item(5, 38).Value2 = 'OKC'

# Row 6: DEN / Michael Porter Jr. (unchanged, re-asserted)
$ws.Cells.Item(6, 1).Value2 = 'DEN'
$ws.Cells.Item(6, 2).Value2 = 'Michael Porter Jr.'
$ws.Cells.Item(6, 3).Value2 = $null
$ws.Cells.Item(6, 4).Value2 = 'F'
$ws.Cells.Item(6, 5).Value2 = $null
$ws.Cells.Item(6, 6).Value2 = 32.2
$ws.Cells.Item(6, 7).Value2 = 30.8
$ws.Cells.Item(6, 8).Value2 = 24.6
$ws.Cells.Item(6, 9).Value2 = 14
$ws.Cells.Item(6, 10).Value2 = 3
$ws.Cells.Item(6, 11).Value2 = 3
$ws.Cells.Item(6, 12).Value2 = 5
$ws.Cells.Item(6, 13).Value2 = 3
$ws.Cells.Item(6, 14).Value2 = 35
$ws.Cells.Item(6, 15).Value2 = 35
$ws.Cells.Item(6, 16).Value2 = 23
$ws.Cells.Item(6, 17).Value2 = 36
$ws.Cells.Item(6, 18).Value2 = 32
$ws.Cells.Item(6, 19).Value2 = 'vs'
$ws.Cells.Item(6, 20).Value2 = 0.9
$ws.Cells.Item(6, 21).Value2 = $null
$ws.Cells.Item(6, 22).Value2 = $null
$ws.Cells.Item(6, 23).Value2 = 'vs'
$ws.Cells.Item(6, 24).Value2 = 'NYK'
$ws.Cells.Item(6, 25).Value2 = '@'
$ws.Cells.Item(6, 26).Value2 = 10
$ws.Cells.Item(6, 27).Value2 = '-'
$ws.Cells.Item(6, 28).Value2 = '-'
$ws.Cells.Item(6, 29).Value2 = '-'
$ws.Cells.Item(6, 30).Value2 = '-'
$ws.Cells.Item(6, 31).Value2 = '@'
$ws.Cells.Item(6, 32).Value2 = 'POR'
$ws.Cells.Item(6, 33).Value2 = 'vs'
$ws.Cells.Item(6, 34).Value2 = 'MEM'
$ws.Cells.Item(6, 35).Value2 = 'vs'
$ws.Cells.Item(6, 36).Value2 = 'PHX'
$ws.Cells.Item(6, 37).Value2 = 'vs'
$ws.Cells.Item(6, 38).Value2 = 'MIN'

# Row 7: CHI / Nikola Vucevic (unchanged, re-asserted)
$ws.Cells.Item(7, 1).Value2 = 'CHI'
$ws.Cells.Item(7, 2).Value2 = 'Nikola Vucevic'
$ws.Cells.Item(7, 3).Value2 = $null
$ws.Cells.Item(7, 4).Value2 = 'C'
$ws.Cells.Item(7, 5).Value2 = $null
$ws.Cells.Item(7, 6).Value2 = 30.2
$ws.Cells.Item(7, 7).Value2 = 30.8
$ws.Cells.Item(7, 8).Value2 = 30.6
$ws.Cells.Item(7, 9).Value2 = 14
$ws.Cells.Item(7, 10).Value2 = 2
$ws.Cells.Item(7, 11).Value2 = 2
$ws.Cells.Item(7, 12).Value2 = 8
$ws.Cells.Item(7, 13).Value2 = 2
$ws.Cells.Item(7, 14).Value2 = 32
$ws.Cells.Item(7, 15).Value2 = 50
$ws.Cells.Item(7, 16).Value2 = 32
$ws.Cells.Item(7, 17).Value2 = 21
$ws.Cells.Item(7, 18).Value2 = 16
$ws.Cells.Item(7, 19).Value2 = '@'
$ws.Cells.Item(7, 20).Value2 = 1.6
$ws.Cells.Item(7, 21).Value2 = $null
$ws.Cells.Item(7, 22).Value2 = $null
$ws.Cells.Item(7, 23).Value2 = '@'
$ws.Cells.Item(7, 24).Value2 = 'HOU'
$ws.Cells.Item(7, 25).Value2 = 'vs'
$ws.Cells.Item(7, 26).Value2 = 33
$ws.Cells.Item(7, 27).Value2 = '-'
$ws.Cells.Item(7, 28).Value2 = '-'
$ws.Cells.Item(7, 29).Value2 = '-'
$ws.Cells.Item(7, 30).Value2 = '-'
$ws.Cells.Item(7, 31).Value2 = 'vs'
$ws.Cells.Item(7, 32).Value2 = 'BOS'
$ws.Cells.Item(7, 33).Value2 = 'vs'
$ws.Cells.Item(7, 34).Value2 = 'WAS'
$ws.Cells.Item(7, 35).Value2 = 'vs'
$ws.Cells.Item(7, 36).Value2 = 'IND'
$ws.Cells.Item(7, 37).Value2 = '@'
$ws.Cells.Item(7, 38).Value2 = 'BKN'

# Row 8: UTA / Lauri Markkanen (status Out -> Day-To-Day, M-1/M-2 refreshed)
$ws.Cells.Item(8, 1).Value2 = 'UTA'
$ws.Cells.Item(8, 2).Value2 = 'Lauri Markkanen'
$ws.Cells.Item(8, 3).Value2 = 'Day-To-Day'
$ws.Cells.Item(8, 4).Value2 = 'F'
$ws.Cells.Item(8, 5).Value2 = 'O'
$ws.Cells.Item(8, 6).Value2 = 28.8
$ws.Cells.Item(8, 7).Value2 = 29.3
$ws.Cells.Item(8, 8).Value2 = 35
$ws.Cells.Item(8, 9).Value2 = 6
$ws.Cells.Item(8, 10).Value2 = 1
$ws.Cells.Item(8, 11).Value2 = 1
$ws.Cells.Item(8, 12).Value2 = 2
$ws.Cells.Item(8, 13).Value2 = 2
$ws.Cells.Item(8, 14).Value2 = '-'
$ws.Cells.Item(8, 15).Value2 = 34
$ws.Cells.Item(8, 16).Value2 = '-'
$ws.Cells.Item(8, 17).Value2 = '-'
$ws.Cells.Item(8, 18).Value2 = '-'
$ws.Cells.Item(8, 19).Value2 = '@'
$ws.Cells.Item(8, 20).Value2 = -2.6
$ws.Cells.Item(8, 21).Value2 = -1
$ws.Cells.Item(8, 22).Value2 = 7
$ws.Cells.Item(8, 23).Value2 = '@'
$ws.Cells.Item(8, 24).Value2 = 'DAL'
$ws.Cells.Item(8, 25).Value2 = 'vs'
$ws.Cells.Item(8, 26).Value2 = 24
$ws.Cells.Item(8, 27).Value2 = '-'
$ws.Cells.Item(8, 28).Value2 = '-'
$ws.Cells.Item(8, 29).Value2 = '-'
$ws.Cells.Item(8, 30).Value2 = '-'
$ws.Cells.Item(8, 31).Value2 = '@'
$ws.Cells.Item(8, 32).Value2 = 'HOU'
$ws.Cells.Item(8, 33).Value2 = 'vs'
$ws.Cells.Item(8, 34).Value2 = 'DAL'
$ws.Cells.Item(8, 35).Value2 = 'vs'
$ws.Cells.Item(8, 36).Value2 = 'SAS'
$ws.Cells.Item(8, 37).Value2 = 'vs'
$ws.Cells.Item(8, 38).Value2 = 'HOU'

# Row 11: BKN / Nic Claxton (status -> Questionable)
$ws.Cells.Item(11, 1).Value2 = 'BKN'
$ws.Cells.Item(11, 2).Value2 = 'Nic Claxton'
$ws.Cells.Item(11, 3).Value2 = 'Questionable'
$ws.Cells.Item(11, 4).Value2 = 'C'
$ws.Cells.Item(11, 5).Value2 = $null
$ws.Cells.Item(11, 6).Value2 = 26
$ws.Cells.Item(11, 7).Value2 = 26.7
$ws.Cells.Item(11, 8).Value2 = 28
$ws.Cells.Item(11, 9).Value2 = 15
$ws.Cells.Item(11, 10).Value2 = 2
$ws.Cells.Item(11, 11).Value2 = 7
$ws.Cells.Item(11, 12).Value2 = 6
$ws.Cells.Item(11, 13).Value2 = 0
$ws.Cells.Item(11, 14).Value2 = 21
$ws.Cells.Item(11, 15).Value2 = 27
$ws.Cells.Item(11, 16).Value2 = 24
$ws.Cells.Item(11, 17).Value2 = 23
$ws.Cells.Item(11, 18).Value2 = 35
$ws.Cells.Item(11, 19).Value2 = '@'
$ws.Cells.Item(11, 20).Value2 = 0
$ws.Cells.Item(11, 21).Value2 = $null
$ws.Cells.Item(11, 22).Value2 = $null
$ws.Cells.Item(11, 23).Value2 = '@'
$ws.Cells.Item(11, 24).Value2 = 'MIL'
$ws.Cells.Item(11, 25).Value2 = '-'
$ws.Cells.Item(11, 26).Value2 = '-'
$ws.Cells.Item(11, 27).Value2 = '-'
$ws.Cells.Item(11, 28).Value2 = '-'
$ws.Cells.Item(11, 29).Value2 = '-'
$ws.Cells.Item(11, 30).Value2 = '-'
$ws.Cells.Item(11, 31).Value2 = '@'
$ws.Cells.Item(11, 32).Value2 = 'NYK'
$ws.Cells.Item(11, 33).Value2 = '@'
$ws.Cells.Item(11, 34).Value2 = 'TOR'
$ws.Cells.Item(11, 35).Value2 = '@'
$ws.Cells.Item(11, 36).Value2 = 'WAS'
$ws.Cells.Item(11, 37).Value2 = 'vs'
$ws.Cells.Item(11, 38).Value2 = 'CHI'

# Row 12: MIL / Khris Middleton (refreshed stats)
$ws.Cells.Item(12, 1).Value2 = 'MIL'
$ws.Cells.Item(12, 2).Value2 = 'Khris Middleton'
$ws.Cells.Item(12, 3).Value2 = $null
$ws.Cells.Item(12, 4).Value2 = 'F'
$ws.Cells.Item(12, 5).Value2 = 'O'
$ws.Cells.Item(12, 6).Value2 = 28.8
$ws.Cells.Item(12, 7).Value2 = 25.9
$ws.Cells.Item(12, 8).Value2 = 24
$ws.Cells.Item(12, 9).Value2 = 2
$ws.Cells.Item(12, 10).Value2 = 0
$ws.Cells.Item(12, 11).Value2 = 0
$ws.Cells.Item(12, 12).Value2 = 1
$ws.Cells.Item(12, 13).Value2 = 1
$ws.Cells.Item(12, 14).Value2 = 41
$ws.Cells.Item(12, 15).Value2 = 36
$ws.Cells.Item(12, 16).Value2 = '-'
$ws.Cells.Item(12, 17).Value2 = '-'
$ws.Cells.Item(12, 18).Value2 = '-'
$ws.Cells.Item(12, 19).Value2 = 'vs'
$ws.Cells.Item(12, 20).Value2 = 0.9
$ws.Cells.Item(12, 21).Value2 = -1
$ws.Cells.Item(12, 22).Value2 = 1
$ws.Cells.Item(12, 23).Value2 = 'vs'
$ws.Cells.Item(12, 24).Value2 = 'BKN'
$ws.Cells.Item(12, 25).Value2 = '@'
$ws.Cells.Item(12, 26).Value2 = 44
$ws.Cells.Item(12, 27).Value2 = '@'
$ws.Cells.Item(12, 28).Value2 = 24
$ws.Cells.Item(12, 29).Value2 = '-'
$ws.Cells.Item(12, 30).Value2 = '-'
$ws.Cells.Item(12, 31).Value2 = 'vs'
$ws.Cells.Item(12, 32).Value2 = 'OKC'
$ws.Cells.Item(12, 33).Value2 = 'vs'
$ws.Cells.Item(12, 34).Value2 = 'LAL'
$ws.Cells.Item(12, 35).Value2 = '@'
$ws.Cells.Item(12, 36).Value2 = 'NOP'
$ws.Cells.Item(12, 37).Value2 = '@'
$ws.Cells.Item(12, 38).Value2 = 'ATL'
